$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlanetColor")

# PLANET_ONE row (row 4): drop trailing colon from the localized labels and
# give French its own dedicated "Planète 1" string (was wrongly sharing the
# PlanetOpacity "Planète 1:" string). German (I4) reused the English label,
# which is preserved, just without the colon now too.
$ws.Range("B4").Value = "Planet 1"
$ws.Range("C4").Value = "행성 1"
$ws.Range("D4").Value = "Planeta 1"
$ws.Range("F4").Value = "Planète 1"
$ws.Range("G4").Value = "Hành tinh thứ nhất"
$ws.Range("I4").Value = "Planet 1"

# PLANET_TWO row (row 5): same fix.
$ws.Range("B5").Value = "Planet 2"
$ws.Range("C5").Value = "행성 2"
$ws.Range("D5").Value = "Planeta 2"
$ws.Range("F5").Value = "Planète 2"
$ws.Range("G5").Value = "Hành tinh thứ hai"
$ws.Range("I5").Value = "Planet 2"
